$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Fill in the annotation scores for columns E:J, rows 2-29 ----
$data = @{
    2  = @(2,2,2,2,2,2)
    3  = @(2,2,2,2,1,2)
    4  = @(1,1,1,2,2,1)
    5  = @(2,2,1,1,1,2)
    6  = @(2,1,1,1,1,2)
    7  = @(2,2,2,2,2,2)
    8  = @(2,1,1,1,1,2)
    9  = @(2,2,1,1,2,2)
    10 = @(2,2,1,1,2,2)
    11 = @(2,1,1,1,1,2)
    12 = @(2,1,1,1,1,2)
    13 = @(1,1,1,1,0,1)
    14 = @(2,1,1,2,2,2)
    15 = @(2,1,1,1,1,2)
    16 = @(2,2,1,2,2,2)
    17 = @(2,1,0,1,2,2)
    18 = @(2,1,1,1,1,2)
    19 = @(2,1,1,1,1,2)
    20 = @(2,2,2,2,1,2)
    21 = @(2,2,1,2,2,2)
    22 = @(2,1,1,1,2,2)
    23 = @(2,1,1,1,0,2)
    24 = @(1,1,1,2,1,1)
    25 = @(2,1,1,1,1,2)
    26 = @(1,1,1,2,1,2)
    27 = @(2,1,1,1,1,2)
    28 = @(1,1,1,1,1,1)
    29 = @(1,1,1,1,1,1)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i  # column E is 5
        $cell = $ws.Cells.Item($r, $col)
        $cell.Value = $vals[$i]
        $cell.Style = "Normal"
    }
}

# Ensure row spans metadata / style consistency: apply the same style (s=2) used by A:D already
# by copying number format/alignment from an existing styled cell (e.g. A2) to E2:J29.
$srcStyleRange = $ws.Range("D2")
$destStyleRange = $ws.Range("E2:J29")
$srcStyleRange.Copy() | Out-Null
$destStyleRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# ---- Update the sheet view: zoom, freeze the header row, and selection ----
$ws.Activate()
$window = $excel.ActiveWindow
# Freezing happens relative to the currently selected cell (row 2 => 1 frozen row).
$ws.Range("A2").Select() | Out-Null
$window.FreezePanes = $true
$window.Zoom = 85
$window.ScrollRow = 20
$ws.Range("E30").Select() | Out-Null
